$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training")

$ws.Range("D2").Value = 0.008880000000000001
$ws.Range("D3").Value = 0.04308
$ws.Range("D4").Value = 0.98027
$ws.Range("D5").Value = 0.00886
$ws.Range("D6").Value = 0.04363
$ws.Range("D7").Value = 0.99717
$ws.Range("D8").Value = 0.01071
$ws.Range("D9").Value = 0.03814
$ws.Range("D10").Value = 1.01244
